# CTMS: Added visit schedule test case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 4) mirroring row 1, plus an extra "Category" column (F)
$ws.Range("A4").Value = "AddVisitScheduleForSite"
$ws.Range("B4").Value = "User"
$ws.Range("C4").Value = "Password"
$ws.Range("D4").Value = "Study Phase"
$ws.Range("E4").Value = "Status"

# New data row (row 5) mirroring row 2, plus the new "Category" value (F5)
$ws.Range("B5").Value = "usersetup"
$ws.Range("C5").Value = "b1f0rcE"
$ws.Range("D5").Value = "Phase II/III"
$ws.Range("E5").Value = "Planning"
$ws.Range("F5").Value = "Qualification Visit"
$ws.Range("F4").Value = "Category"

# Match formatting used by the existing header/credential rows
$ws.Range("B4:F4").Font.Bold = $true
$ws.Range("B5:C5").Style = "Hyperlink"

# Update the active selection to the newly added cell
$null = $ws.Range("F5").Select()
